$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 847
$ws.Range("F3").Value = 1419
$ws.Range("F4").Value = 1002
$ws.Range("F5").Value = 490
$ws.Range("F6").Value = 203
$ws.Range("F7").Value = 642
$ws.Range("F8").Value = 210
$ws.Range("F10").Value = 55
$ws.Range("F12").Value = 130
$ws.Range("F13").Value = 1688
$ws.Range("F14").Value = 296
$ws.Range("F16").Value = 481
$ws.Range("F17").Value = 81
$ws.Range("F18").Value = 401
$ws.Range("F21").Value = 640
$ws.Range("F22").Value = 38
$ws.Range("F23").Value = 224
$ws.Range("F24").Value = 942
$ws.Range("F26").Value = 1495
$ws.Range("F27").Value = 216

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 118
$ws.Range("F3").Value = 32
$ws.Range("F4").Value = 655
$ws.Range("F5").Value = 205
$ws.Range("F6").Value = 16
$ws.Range("F7").Value = 277

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 847
$ws.Range("F4").Value = 1419
$ws.Range("F5").Value = 1002
$ws.Range("F6").Value = 118
$ws.Range("F7").Value = 32
$ws.Range("F8").Value = 490
$ws.Range("F9").Value = 203
$ws.Range("F10").Value = 642
$ws.Range("F11").Value = 655
$ws.Range("F12").Value = 210
$ws.Range("F14").Value = 55
$ws.Range("F16").Value = 130
$ws.Range("F17").Value = 1688
$ws.Range("F18").Value = 205
$ws.Range("F19").Value = 296
$ws.Range("F21").Value = 481
$ws.Range("F22").Value = 81
$ws.Range("F23").Value = 401
$ws.Range("F24").Value = 16
$ws.Range("F27").Value = 277
$ws.Range("F29").Value = 640
$ws.Range("F34").Value = 38
$ws.Range("F35").Value = 224
$ws.Range("F36").Value = 942
$ws.Range("F38").Value = 1495
$ws.Range("F39").Value = 216
